# "Added Week 15 simulations"
# Appends one more week's worth of simulated per-game numbers to the
# long space-separated simulation strings on YDS / ST, and updates the
# season-to-date aggregate totals on OFF / DEF / ST / TURNS / PEN that
# those new simulated games feed into.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append newly-simulated per-game yardage numbers
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value2 + " 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4"
$ws.Range("B3").Value = $ws.Range("B3").Value2 + " 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11"
$ws.Range("C2").Value = $ws.Range("C2").Value2 + " 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3"
$ws.Range("C3").Value = $ws.Range("C3").Value2 + " -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23"

# ---------------------------------------------------------------
# OFF sheet: updated season totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 175
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 44
$ws.Range("G2").Value = 51
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 30
$ws.Range("L2").Value = 245
$ws.Range("M2").Value = 145
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 11
$ws.Range("Q2").Value = 438

$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 132
$ws.Range("E3").Value = 31
$ws.Range("F3").Value = 97
$ws.Range("I3").Value = 57
$ws.Range("J3").Value = 34
$ws.Range("N3").Value = 13

# ---------------------------------------------------------------
# DEF sheet: updated season totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 163
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 52
$ws.Range("G2").Value = 42
$ws.Range("J2").Value = 26
$ws.Range("L2").Value = 211
$ws.Range("M2").Value = 133
$ws.Range("O2").Value = 19
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = 430

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 112
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 29
$ws.Range("F3").Value = 79
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 38

# ---------------------------------------------------------------
# ST sheet: updated special-teams totals plus simulation strings
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 61
$ws.Range("D2").Value = 56
$ws.Range("F2").Value = 134
$ws.Range("G2").Value = 127

$ws.Range("B3").Value = 40

$ws.Range("B4").Value = $ws.Range("B4").Value2 + " 63"
$ws.Range("B5").Value = $ws.Range("B5").Value2 + " 38"
$ws.Range("B6").Value = $ws.Range("B6").Value2 + " 26 28 0 25"

$ws.Range("D3").Value = $ws.Range("D3").Value2 + " 54 33"
$ws.Range("D4").Value = $ws.Range("D4").Value2 + " 0 0"
$ws.Range("D5").Value = $ws.Range("D5").Value2 + " 0 0"

# ---------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")

$ws.Range("B2").Value = 10
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 11
$ws.Range("E3").Value = 10

# ---------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")

$ws.Range("B2").Value = 19
